$wb = $excel.ActiveWorkbook

# --- 1. Append rows 9-11 to the existing "Peerapat" sheet ---
$peerapat = $wb.Worksheets.Item("Peerapat")

$peerapat.Range("A9").Value = "04/10/2021 21:45"
$peerapat.Range("B9").Value = 50
$peerapat.Range("C9").Value = 170
$peerapat.Range("D9").Value = 17.3
$peerapat.Range("E9").Value = "ผอม"

$peerapat.Range("A10").Value = "04/10/2021 21:45"
$peerapat.Range("B10").Value = 80
$peerapat.Range("C10").Value = 170
$peerapat.Range("D10").Value = 27.68
$peerapat.Range("E10").Value = "อ้วน"

$peerapat.Range("A11").Value = "04/10/2021 21:46"
$peerapat.Range("B11").Value = 100
$peerapat.Range("C11").Value = 170
$peerapat.Range("D11").Value = 34.6
$peerapat.Range("E11").Value = "อ้วนมาก"

# --- 2. Add a new "Chanakarn" sheet after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$chanakarn = $wb.Worksheets.Add($null, $lastSheet)
$chanakarn.Name = "Chanakarn"

$chanakarn.Range("A1").Value = "Chanakarn"

$chanakarn.Range("A2").Value = "เวลา"
$chanakarn.Range("B2").Value = "น้ำหนัก(กก.)"
$chanakarn.Range("C2").Value = "ส่วนสูง(ซม.)"
$chanakarn.Range("D2").Value = "BMI"
$chanakarn.Range("E2").Value = "เกณฑ์"

$chanakarn.Range("A3").Value = "04/10/2021 21:46"
$chanakarn.Range("B3").Value = 100
$chanakarn.Range("C3").Value = 170
$chanakarn.Range("D3").Value = 34.6
$chanakarn.Range("E3").Value = "อ้วนมาก"

# --- 3. Add a new "Jedilok" sheet after "Chanakarn" ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$jedilok = $wb.Worksheets.Add($null, $lastSheet2)
$jedilok.Name = "Jedilok"

$jedilok.Range("A1").Value = "Jedilok"

$jedilok.Range("A2").Value = "เวลา"
$jedilok.Range("B2").Value = "น้ำหนัก(กก.)"
$jedilok.Range("C2").Value = "ส่วนสูง(ซม.)"
$jedilok.Range("D2").Value = "BMI"
$jedilok.Range("E2").Value = "เกณฑ์"
